# g3.9 - remoção de ano no nome da aba para evitar futuros erros

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet, dropping the trailing " 2023"
$ws.Name = "g3.11a"

# Add the "Ano" column header (D1), copying the style of the other headers
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("D1").Value = "Ano"

# Fill D2:D9 with the year value 2023 (numeric)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value = 2023
}
